# "ispunjen 11. zahtev i ulepsan ostatak koda"
#
# Sheet2 had a merged range G24:G26 (value 6, centered) used for the
# SUBTOTAL-of-G demo. The edit un-merges that range and instead puts a
# distinct value in each of the three G cells (G25=6, G26=7, G24 stays
# blank), clearing out the old E24/F24/G24 numbers that fed the old
# layout. The SUBTOTAL formulas in row 28 recalculate on their own.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Break up the old G24:G26 merged cell.
$ws.Range("G24:G26").UnMerge()

# Row 24 no longer carries any of the old demo numbers.
$ws.Range("E24").ClearContents()
$ws.Range("F24").ClearContents()
$ws.Range("G24").ClearContents()

# New per-row values for the previously-merged G column cells.
$ws.Range("G25").Value = 6
$ws.Range("G26").Value = 7

# Keep the (now distinct) G24:G26 formatting as general/horizontal-default
# alignment, matching the style split that unmerging produced.
$ws.Range("G24:G26").HorizontalAlignment = 1

# Match the author's final selection state.
$ws.Range("M22").Select()

$wb.Application.Calculate()
